$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to remain plain text, matching
# the original inline-string cell type, so numeric-looking values (e.g.
# "1.001", "0.06469") are not coerced into numbers by Excel's type inference.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.502.92"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.885.18"
$ws.Range("E3").Value = "  +0.79%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
$ws.Range("D5").Value = "244.18"
$ws.Range("E5").Value = "  -1.28%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4697"
$ws.Range("E7").Value = "  -0.84%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2889"
$ws.Range("E8").Value = "  -0.61%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06469"
$ws.Range("E9").Value = "  +0.02%  "

# Row 10 - Solana
$ws.Range("D10").Value = "22.24"
$ws.Range("E10").Value = "  +0.83%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07752"
$ws.Range("E11").Value = "  +0.49%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.883.26"
$ws.Range("E12").Value = "  +0.67%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "0.7296"
$ws.Range("E13").Value = "  -1.19%  "

# Row 14 - Litecoin
$ws.Range("D14").Value = "94.85"
$ws.Range("E14").Value = "  -1.66%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "5.161"
$ws.Range("E15").Value = "  +0.41%  "

# Row 16 - BitcoinCash
$ws.Range("D16").Value = "281.92"
$ws.Range("E16").Value = "  +3.46%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "30.509.19"
$ws.Range("E17").Value = "  -0.14%  "

# Row 18 - Avalanche
$ws.Range("D18").Value = "12.96"
$ws.Range("E18").Value = "  -2.49%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.000007455"
$ws.Range("E20").Value = "  -0.52%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.126.66"
$ws.Range("E21").Value = "  +0.58%  "

# Row 22 - BinanceUSD
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.26%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.235"
$ws.Range("E23").Value = "  +0.02%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "6.226"
$ws.Range("E24").Value = "  +0.88%  "

# Row 25 - was Cosmos, becomes Monero
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "163.63"
$ws.Range("E25").Value = "  +0.19%  "

# Row 26 - was Monero, becomes Cosmos
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.075"
$ws.Range("E26").Value = "  -1.42%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").Value = "  +0.20%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "1.886"
$ws.Range("E28").Value = "  -1.29%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "1.333"
$ws.Range("E29").Value = "  -0.85%  "

# Row 30 - Stellar
$ws.Range("D30").Value = "0.09697"
$ws.Range("E30").Value = "  -2.72%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "1.470"
$ws.Range("E31").Value = "  -2.43%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "4.256"
$ws.Range("E32").Value = "  -0.60%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "4.095"
$ws.Range("E33").Value = "  -0.12%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.04854"
$ws.Range("E34").Value = "  +1.54%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "1.122"
$ws.Range("E35").Value = "  +0.45%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.6888"
$ws.Range("E36").Value = "  -0.86%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "2.707"
$ws.Range("E37").Value = "  -0.39%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01891"
$ws.Range("E38").Value = "  +2.52%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "2.820"
$ws.Range("E39").Value = "  +2.51%  "

# Row 40 - Aave
$ws.Range("D40").Value = "75.10"
$ws.Range("E40").Value = "  +2.66%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "6.149"
$ws.Range("E41").Value = "  -0.51%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "1.994"
$ws.Range("E42").Value = "  +1.44%  "

# Row 43 - TheSandbox
$ws.Range("D43").Value = "0.4234"
$ws.Range("E43").Value = "  +1.58%  "

# Row 44 - PaxDollar
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.06%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").Value = "0.8215"
$ws.Range("E45").Value = "  -1.36%  "

# Row 46 - Quant
$ws.Range("D46").Value = "100.83"
$ws.Range("E46").Value = "  -1.68%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "9.516"
$ws.Range("E47").Value = "  +2.30%  "

# Row 48 - Elrond
$ws.Range("D48").Value = "35.19"
$ws.Range("E48").Value = "  -0.53%  "

# Row 49 - Aptos
$ws.Range("D49").Value = "6.950"
$ws.Range("E49").Value = "  +0.06%  "

# Row 50 - Maker
$ws.Range("D50").Value = "910.13"
$ws.Range("E50").Value = "  -0.98%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  +1.83%  "
